# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / LevePrice* / LeveProfit* derived columns
# (H..N) for the rows whose underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 290.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H88").Value = 10004
$ws.Range("J88").Value = 10004
$ws.Range("L88").Value = 10004
$ws.Range("N88").Value = -10816

$ws.Range("H91").Value = 10004
$ws.Range("J91").Value = 10004
$ws.Range("L91").Value = 10004
$ws.Range("N91").Value = -12812

$ws.Range("H132").Value = 2071.6206
$ws.Range("I132").Value = 2113.963
$ws.Range("K132").Value = 6341.889000000001
$ws.Range("M132").Value = -3811.889000000001

$ws.Range("H137").Value = 2699.158
$ws.Range("J137").Value = 4840.4
$ws.Range("L137").Value = 14521.2
$ws.Range("N137").Value = -19621.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8327.772000000001
$ws.Range("I32").Value = 6013.5312
$ws.Range("J32").Value = 14499.083
$ws.Range("K32").Value = 6013.5312
$ws.Range("L32").Value = 14499.083
$ws.Range("M32").Value = -5726.5312
$ws.Range("N32").Value = -15073.083

$ws.Range("H61").Value = 2179.647
$ws.Range("I61").Value = 2068.4285
$ws.Range("J61").Value = 2698.6667
$ws.Range("K61").Value = 2068.4285
$ws.Range("L61").Value = 2698.6667
$ws.Range("M61").Value = -1856.4285
$ws.Range("N61").Value = -3122.6667

$ws.Range("H132").Value = 1670.0278
$ws.Range("I132").Value = 1579
$ws.Range("K132").Value = 4737
$ws.Range("M132").Value = -2207

$ws.Range("H136").Value = 2179.647
$ws.Range("I136").Value = 2068.4285
$ws.Range("J136").Value = 2698.6667
$ws.Range("K136").Value = 6205.2855
$ws.Range("L136").Value = 8096.000100000001
$ws.Range("M136").Value = -3655.2855
$ws.Range("N136").Value = -13196.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1774.5
$ws.Range("I86").Value = 1718.6923
$ws.Range("K86").Value = 1718.6923
$ws.Range("M86").Value = -595.6922999999999

$ws.Range("H89").Value = 1774.5
$ws.Range("I89").Value = 1718.6923
$ws.Range("K89").Value = 8593.461499999999
$ws.Range("M89").Value = -2977.461499999999

$ws.Range("H134").Value = 3111.3572
$ws.Range("J134").Value = 2961.2856
$ws.Range("L134").Value = 8883.856800000001
$ws.Range("N134").Value = -13953.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5641.1055
$ws.Range("I31").Value = 4619.6665
$ws.Range("K31").Value = 4619.6665
$ws.Range("M31").Value = -4324.6665

$ws.Range("H34").Value = 5641.1055
$ws.Range("I34").Value = 4619.6665
$ws.Range("K34").Value = 4619.6665
$ws.Range("M34").Value = -4417.6665

$ws.Range("H58").Value = 2865.5173
$ws.Range("J58").Value = 4034.5881
$ws.Range("L58").Value = 4034.5881
$ws.Range("N58").Value = -4440.5881

$ws.Range("H99").Value = 12400.954
$ws.Range("I99").Value = 7528
$ws.Range("K99").Value = 7528
$ws.Range("M99").Value = -6030

$ws.Range("H107").Value = 1011
$ws.Range("I107").Value = 1011
$ws.Range("K107").Value = 1011
$ws.Range("M107").Value = 909

$ws.Range("H126").Value = 12400.954
$ws.Range("I126").Value = 7528
$ws.Range("K126").Value = 22584
$ws.Range("M126").Value = -20114

$ws.Range("H132").Value = 2174.2896
$ws.Range("I132").Value = 1962.8108
$ws.Range("K132").Value = 5888.4324
$ws.Range("M132").Value = -3358.4324

$ws.Range("H134").Value = 2686
$ws.Range("I134").Value = 1523.4445
$ws.Range("K134").Value = 4570.333500000001
$ws.Range("M134").Value = -2035.333500000001

$ws.Range("H136").Value = 2865.5173
$ws.Range("J136").Value = 4034.5881
$ws.Range("L136").Value = 12103.7643
$ws.Range("N136").Value = -17203.7643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1451.1538
$ws.Range("I33").Value = 1660.7142
$ws.Range("J33").Value = 1206.6666
$ws.Range("K33").Value = 9964.285199999998
$ws.Range("L33").Value = 7239.9996
$ws.Range("M33").Value = -9681.285199999998
$ws.Range("N33").Value = -7805.9996

$ws.Range("H115").Value = 1374.5
$ws.Range("J115").Value = 1250
$ws.Range("L115").Value = 3750
$ws.Range("N115").Value = -6100

$ws.Range("H122").Value = 302.33334
$ws.Range("I122").Value = 351.5
$ws.Range("K122").Value = 3163.5
$ws.Range("M122").Value = -713.5

$ws.Range("H139").Value = 7249.75
$ws.Range("I139").Value = 4000
$ws.Range("K139").Value = 12000
$ws.Range("M139").Value = -6860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 102865.4
$ws.Range("I122").Value = 2786.6
$ws.Range("J122").Value = 202944.2
$ws.Range("K122").Value = 8359.799999999999
$ws.Range("L122").Value = 608832.6000000001
$ws.Range("M122").Value = -5909.799999999999
$ws.Range("N122").Value = -613732.6000000001

$ws.Range("H123").Value = 56333.668
$ws.Range("J123").Value = 56333.668
$ws.Range("L123").Value = 56333.668
$ws.Range("N123").Value = -61233.668

$ws.Range("H132").Value = 3206.4443
$ws.Range("I132").Value = 2330.182
$ws.Range("J132").Value = 4583.4287
$ws.Range("K132").Value = 6990.545999999999
$ws.Range("L132").Value = 13750.2861
$ws.Range("M132").Value = -4460.545999999999
$ws.Range("N132").Value = -18810.2861

$ws.Range("H134").Value = 87282
$ws.Range("J134").Value = 87282
$ws.Range("L134").Value = 261846
$ws.Range("N134").Value = -266916

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5833.5
$ws.Range("I122").Value = 6410.4443
$ws.Range("J122").Value = 5091.7144
$ws.Range("K122").Value = 19231.3329
$ws.Range("L122").Value = 15275.1432
$ws.Range("M122").Value = -16781.3329
$ws.Range("N122").Value = -20175.1432

$ws.Range("H136").Value = 2942.04
$ws.Range("I136").Value = 2978.8096
$ws.Range("J136").Value = 2749
$ws.Range("K136").Value = 8936.4288
$ws.Range("L136").Value = 8247
$ws.Range("M136").Value = -6386.4288
$ws.Range("N136").Value = -13347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 966.1539
$ws.Range("I113").Value = 932.7273
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 2798.1819
$ws.Range("L113").Value = 3450
$ws.Range("M113").Value = -628.1819
$ws.Range("N113").Value = -7790

$ws.Range("H122").Value = 3792.875
$ws.Range("I122").Value = 3792.875
$ws.Range("K122").Value = 11378.625
$ws.Range("M122").Value = -8928.625

$ws.Range("H136").Value = 1550.1786
$ws.Range("I136").Value = 1577.1538
$ws.Range("K136").Value = 4731.4614
$ws.Range("M136").Value = -2181.4614
